$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = 2
$ws.Range("F6").Value = -2
